$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer text (A13)
$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-30 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-10
$ws.Range("D2").Value = 0.1001829102659257
$ws.Range("E2").Value = 0.01637370575487584

$ws.Range("D3").Value = 0.1042563033711244
$ws.Range("E3").Value = 0.001334519572953496

$ws.Range("D4").Value = 0.1182897864653878
$ws.Range("E4").Value = -0.005730189914865846

$ws.Range("D5").Value = 0.137738718112805
$ws.Range("E5").Value = 0.002893445716442233

$ws.Range("D6").Value = 0.1365193163399668
$ws.Range("E6").Value = -0.008729250143102552

$ws.Range("D7").Value = 0.1457961472136296
$ws.Range("E7").Value = 0.005774689243417841

$ws.Range("D8").Value = 0.1285094537613482
$ws.Range("E8").Value = 0.0006080875646092831

$ws.Range("D9").Value = 0.1287073644698125
$ws.Range("E9").Value = -0.001662119622245495

$ws.Range("E10").Value = 0.001014648280393882
